$d = $word.ActiveDocument

# Locate the "SMARTREWARDS FAQ'S" heading paragraph and the last paragraph of
# the FAQ block ("...1000 addresses will get paid.") by scanning paragraph
# text, so the deletion is anchored to content rather than brittle fixed
# indices.
$count = $d.Paragraphs.Count
$titleIdx = -1
$lastFaqIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*SMARTREWARDS FAQ*") {
        $titleIdx = $i
    }
    if ($t -like "*1000 addresses will get paid.*") {
        $lastFaqIdx = $i
    }
}

if ($titleIdx -gt 0 -and $lastFaqIdx -gt $titleIdx) {
    $startPara = $d.Paragraphs.Item($titleIdx + 1)
    $endPara = $d.Paragraphs.Item($lastFaqIdx)
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
